$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.301569700241089
$ws.Range("B1").Value = 2.956575870513916
$ws.Range("C1").Value = 2.327858924865723
$ws.Range("D1").Value = 2.166896104812622
$ws.Range("E1").Value = 1.838370680809021
